# Rename Chert_type category labels from "A"/"B" to "Coarser"/"Finer"
# across the relevant worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "Chert type": column A holds Chert_type values in rows 2-3
$wsChert = $wb.Worksheets.Item("Chert type")
$wsChert.Range("A2").Value = "Coarser"
$wsChert.Range("A3").Value = "Finer"

# Sheet "Chert+Bamboo": column A holds Chert_type values in rows 2-5
$wsCombo = $wb.Worksheets.Item("Chert+Bamboo")
$wsCombo.Range("A2").Value = "Coarser"
$wsCombo.Range("A3").Value = "Coarser"
$wsCombo.Range("A4").Value = "Finer"
$wsCombo.Range("A5").Value = "Finer"
